$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.072934746742249
$ws.Range("B1").Value = 1.734713196754456
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.91253650188446
$ws.Range("E1").Value = 1.156595826148987
